$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at position 134, pushing the existing row 134
# (and everything below it) down by one row.
$ws.Rows.Item(134).Insert()

# Populate the new row 134 with this week's data point (same market /
# category as the surrounding rows, new date + new price figures).
$ws.Cells.Item(134, 1).Value = 4
$ws.Cells.Item(134, 2).Value = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(134, 3).Value = "Los Lagos"
$ws.Cells.Item(134, 4).Value = 44588
$ws.Cells.Item(134, 5).Value = 10
$ws.Cells.Item(134, 6).Value = 100112039
$ws.Cells.Item(134, 7).Value = "Ciboulette"
$ws.Cells.Item(134, 8).Value = "Sin especificar"
$ws.Cells.Item(134, 9).Value = "Primera"
$ws.Cells.Item(134, 10).Value = 80
$ws.Cells.Item(134, 11).Value = 2500
$ws.Cells.Item(134, 12).Value = 3000
$ws.Cells.Item(134, 13).Value = 2750
$ws.Cells.Item(134, 14).Value = "`$/docena de atados"
$ws.Cells.Item(134, 15).Value = "Región Metropolitana"
$ws.Cells.Item(134, 16).Value = 917
$ws.Cells.Item(134, 17).Value = 3
$ws.Cells.Item(134, 18).Value = "Hortaliza"
